# PTW-TimeSheet 2 (1).xlsx - "28-04-2022" sheet corrections.
# The timesheet entries for rows 47-56 (the first shift/block of the day)
# were re-entered: task names, task types and start/end times changed,
# and a previously-blank trailing row (56) was filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("28-04-2022")

# Row 47: Team Meeting (09:00-09:30) -> Reviewed estimation sheet
$ws.Range("B47").Value = "Reviewed estimation sheet"
$ws.Range("C47").Value = "Project"
$ws.Range("D47").Value = 0.375
$ws.Range("E47").Value = 0.39583333333333331

# Row 48: Soft Skill / Non Project -> Worked on HTML layout for login page / Project
$ws.Range("B48").Value = "Worked on HTML layout for login page"
$ws.Range("C48").Value = "Project"
$ws.Range("D48").Value = 0.39930555555555558
$ws.Range("E48").Value = 0.4375

# Row 49: Customer Meeting / Project -> Morning break / Lunch and Break
$ws.Range("B49").Value = "Morning break"
$ws.Range("C49").Value = "Lunch and Break"
$ws.Range("D49").Value = 0.4375
$ws.Range("E49").Value = 0.45833333333333331

# Row 50: Soft Skill / Non Project -> Reviewed the HTML layout / Project
$ws.Range("B50").Value = "Reviewed the HTML layout"
$ws.Range("C50").Value = "Project"
$ws.Range("D50").Value = 0.46527777777777773
$ws.Range("E50").Value = 0.54166666666666663

# Row 51: Lunch / Lunch and Break -> Customer meeting / Project
$ws.Range("B51").Value = "Customer meeting"
$ws.Range("C51").Value = "Project"
$ws.Range("D51").Value = 0.54166666666666663
$ws.Range("E51").Value = 0.57291666666666663
$ws.Range("F51").Formula = "=E51-D51"

# Row 52: Morning Break / Lunch and Break -> Lunch break / Lunch and Break
$ws.Range("B52").Value = "Lunch break"
$ws.Range("C52").Value = "Lunch and Break"
$ws.Range("D52").Value = 0.57291666666666663
$ws.Range("E52").Value = 0.59375
$ws.Range("F52").Formula = "=E52-D52"

# Row 53: Evening Break / Lunch and Break -> Team meeting / Project
$ws.Range("B53").Value = "Team meeting"
$ws.Range("C53").Value = "Project"
$ws.Range("D53").Value = 0.60416666666666663
$ws.Range("E53").Value = 0.63541666666666663

# Row 54: Learned Angular <Topics> / Project -> Exploration on typescript topics / Exploration
$ws.Range("B54").Value = "Exploration on typescript topics "
$ws.Range("C54").Value = "Exploration "
$ws.Range("D54").Value = 0.64583333333333337
$ws.Range("E54").Value = 0.70833333333333337

# Row 55: Working on HTML layout -> Modified changes in Login page(HTML)
$ws.Range("B55").Value = "Modified changes in Login page(HTML)"
$ws.Range("C55").Value = "Project"
$ws.Range("D55").Value = 0.70833333333333337
$ws.Range("E55").Value = 0.72916666666666663

# Row 56: previously blank -> Worked on HTML layout for home page / Project
$ws.Range("B56").Value = "Worked on HTML layout for home page"
$ws.Range("C56").Value = "Project"
$ws.Range("D56").Value = 0.875
$ws.Range("E56").Value = 0.9375

# Update the window's visible scroll/selection for this sheet.
$null = $ws.Range("D52").Select()
